# "add 'UBNT_M2' automation script" -- extend the testData header/value
# rows in the UBNT test workbook with one more tab definition column
# (tabName3 = "U"), appended right after the existing gatewayIP column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U2").Value = "tabName3"
$ws.Range("U3").Value = "U"

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("E3").Select() | Out-Null

Write-Output "added tabName3/U column (U2:U3) to testData sheet"
